$wb = $excel.ActiveWorkbook
$wb.Worksheets.Item("Sheet4").Delete()
$wb.Worksheets.Item("Sheet2").Delete()
$wb.Worksheets.Item("Sheet3").Delete()
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("A27:S27").Interior.Color = 255
